# Batch_Szenarios.xlsx edit — normalize ParBänder (col B) to a constant
# batch size of 4 and raise ParOraclePriceLimit (col K) from 2000 to 10000
# for every scenario row, then leave the cursor on the last touched column
# (matches the "#14 #13 Excel and OnePath" commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Batch_Szenarios")
$ws.Activate()

# ParBänder: every scenario row now uses the same batch value.
$ws.Range("B2:B48").Value = 4

# ParOraclePriceLimit: raised 5x across all scenario rows.
$ws.Range("K2:K48").Value = 10000

# Reflect the author's final cursor position on save.
$ws.Range("J3").Select()
